$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feature rows (12-16) in column A, added first so their shared
# strings land before the later "Marks"/"0-6" strings.
$ws.Range("A12").Value = "Amount of checking for pre- and post- conditions of methods"
$ws.Range("A13").Value = "Provide doctests"
$ws.Range("A14").Value = "Provide unittests"
$ws.Range("A15").Value = "Pretty print, i.e., displaying data in chart/ diagram, etc."
$ws.Range("A16").Value = "Can save and read data from a database, e.g., SQLite, MySQL and MongoDB"

# New "Marks" column header
$ws.Range("D1").Value = "Marks"

# All data rows (2-16) get a "0-6" marks value in column D
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 4).Value = "0-6"
}

# Update the active selection to mirror the final state in the diff
$ws.Range("D16").Select()
